$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C ("Qty executed upto date") numeric updates ---
$ws.Range("C8").Value  = 51
$ws.Range("C9").Value  = 14
$ws.Range("C10").Value = 17
$ws.Range("C11").Value = 26
$ws.Range("C12").Value = 13
$ws.Range("C13").Value = 53
$ws.Range("C14").Value = 62
$ws.Range("C15").Value = 51
$ws.Range("C16").Value = 62
$ws.Range("C17").Value = 12

# --- Column G ("Upto date Amount") text-formatted amounts (kept as text, like the source) ---
# Force each target cell to a Text number format first so the numeric-looking
# string is preserved verbatim (matching the original "###.00" text cells)
# instead of being auto-coerced into a real number by Excel.
$ws.Range("G9").NumberFormat  = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("H21").NumberFormat = "@"

$ws.Range("G9").Value  = "3584.00"
$ws.Range("G10").Value = "8024.00"
$ws.Range("G11").Value = "17212.00"
$ws.Range("G13").Value = "7208.00"
$ws.Range("G14").Value = "1426.00"

# --- Grand-total rows (19 & 21), both G and H columns ---
$ws.Range("G19").Value = "37454.00"
$ws.Range("H19").Value = "37454.00"
$ws.Range("G21").Value = "37454.00"
$ws.Range("H21").Value = "37454.00"
